# Atualização automática via cronjob
# Refreshes the "vendas atipicas" table (A2:H10) with the latest extraction.
# The sheet keeps a rolling window of the most recent atypical-sales rows:
# the oldest day (2025-04-23) is dropped and a new day (2025-05-07) is appended,
# and every row's metrics are refreshed to the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for columns A, C, D, F, G, H (row => values).
# B (Dia) and E (id_produto) are handled separately below because they look
# like dates / numbers to Excel and would otherwise be auto-converted away
# from plain text (losing the "yyyy-mm-dd" layout / leading zeros).
$rows = @(
    @{ Row=2;  A=1; C=150; D='JURUA ESTALEIROS E NAVEGACAO LTDA';     F='VASSOURA PIACAVA 20 FUROS';                               G=174;  H=$false },
    @{ Row=3;  A=5; C=300; D='MUSASHI DA AMAZONIA LTDA';              F='SACO DE LIXO 200L COMUM PACOTINHO C/5 UND';               G=616;  H=$false },
    @{ Row=4;  A=0; C=250; D='MAP SERVICOS DE CONSERVACAO - EIRELI';  F='PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM';       G=-15;  H=$false },
    @{ Row=5;  A=2; C=60;  D='MAP SERVICOS DE CONSERVACAO - EIRELI';  F='INSETICIDA BUZZOFF AEROSSOL 300ML';                       G=2;    H=$true  },
    @{ Row=6;  A=6; C=70;  D='MAP SERVICOS DE CONSERVACAO - EIRELI';  F='DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO';   G=350;  H=$true  },
    @{ Row=7;  A=7; C=20;  D='LUCAS CLIENTE NOVO';                    F='AZULIM LIMPA CERAMICAS E AZULEJOS LAVANDA 5L 1:15 START'; G=0;    H=$true  },
    @{ Row=8;  A=8; C=250; D='MAP SERVICOS DE CONSERVACAO - EIRELI';  F='PEDRA SANITARIA NAFT PLUS FLORAL 25G';                    G=140;  H=$false },
    @{ Row=9;  A=3; C=250; D='RH MULTI SERVICOS ADMINISTRATIVOS S.A'; F='LUVAS DESCARTAVEIS C/ 100 UND';                           G=1055; H=$false },
    @{ Row=10; A=4; C=30;  D='V V REFEICOES LTDA';                    F='CABO DE ALUMINIO NOBRE 140 CM COM PONTEIRA';              G=185;  H=$true  }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("F$($r.Row)").Value = $r.F
    $ws.Range("G$($r.Row)").Value = $r.G
    $ws.Range("H$($r.Row)").Value = $r.H
}

# Helper: write a text value into a set of cells without Excel reinterpreting
# it as a date/number (which would happen with a plain .Value assignment).
# We stage the literal text in an unused cell via a formula (whose result is
# always text), copy that cell onto each target (copy preserves the text
# storage type), then clear the helper cell.
function Set-TextValue($ws, [string]$text, [string[]]$targets) {
    $helper = $ws.Range("Z1")
    $escaped = $text.Replace('"', '""')
    $helper.Formula = "=""$escaped"""
    foreach ($target in $targets) {
        $helper.Copy($ws.Range($target))
    }
    $helper.Clear()
}

# Column B ("Dia") - text dates.
Set-TextValue $ws "2025-04-24" @("B2", "B3")
Set-TextValue $ws "2025-04-28" @("B4", "B5", "B6", "B7", "B8")
Set-TextValue $ws "2025-04-30" @("B9")
Set-TextValue $ws "2025-05-07" @("B10")

# Column E ("id_produto") - zero-padded codes.
Set-TextValue $ws "000088" @("E2")
Set-TextValue $ws "000842" @("E3")
Set-TextValue $ws "000098" @("E4")
Set-TextValue $ws "000347" @("E5")
Set-TextValue $ws "000349" @("E6")
Set-TextValue $ws "000158" @("E7")
Set-TextValue $ws "000779" @("E8")
Set-TextValue $ws "000041" @("E9")
Set-TextValue $ws "000877" @("E10")
